$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a <w:lastRenderedPageBreak/> marker into the run that begins the
#    "Defendant was served ..." sentence (purely structural; no visible text
#    change). We do this via InsertXML on a zero-length range placed right
#    before that text, using a matching rPr so it renders/merges cleanly
#    with the existing run's formatting.
# ---------------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("Defendant was served at least 14 days before the scheduled hearing", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$insPoint = $d.Range($findRng.Start, $findRng.Start)
$insPoint.InsertXML("<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/><w:sz w:val='22'/><w:szCs w:val='22'/></w:rPr><w:lastRenderedPageBreak/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

# ---------------------------------------------------------------------------
# 2) Rework the merge-field conditional logic that used to gate on
#    `uniform_spousal_support_order_foc10b.enabled`; it is replaced
#    everywhere by the combination of `user_wants_post_divorce_support`
#    and `spousal_support_provisions`.
#
#    Several of these {% ... %} tags are duplicated verbatim in two spots
#    in the document (once in the "served with" paragraph, once in the
#    "I REQUEST THAT" paragraph) and always receive the identical
#    replacement text, so wdReplaceAll (2) is used to catch every
#    occurrence in one pass.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute( `
    "{% if (not there_are_marital_children) and (not uniform_spousal_support_order_foc10b.enabled) %}", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    '{% if (not there_are_marital_children) and ((not user_wants_post_divorce_support)and (spousal_support_provisions != "user_agrees_to_pay_spousal_support")) %}', `
    2)

$d.Content.Find.Execute( `
    "{% elif there_are_marital_children and (not uniform_spousal_support_order_foc10b.enabled) %}", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    '{% elif there_are_marital_children and ((not user_wants_post_divorce_support) and (spousal_support_provisions != "user_agrees_to_pay_spousal_support")) %}', `
    2)

$d.Content.Find.Execute( `
    "{% elif (not there_are_marital_children) and uniform_spousal_support_order_foc10b.enabled %}", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    '{% elif (not there_are_marital_children) and (user_wants_post_divorce_support or (spousal_support_provisions == "user_agrees_to_pay_spousal_support")) %}', `
    2)

$d.Content.Find.Execute( `
    "{% elif there_are_marital_children and uniform_spousal_support_order_foc10b.enabled %}", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    '{% elif there_are_marital_children and (user_wants_post_divorce_support or (spousal_support_provisions == "user_agrees_to_pay_spousal_support")) %}', `
    2)

$d.Content.Find.Execute( `
    "{% if there_are_marital_children and uniform_spousal_support_order_foc10b.enabled %}", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    '{% if there_are_marital_children and (user_wants_post_divorce_support or (spousal_support_provisions == "user_agrees_to_pay_spousal_support")) %}', `
    2)
